# Apply the metrics_4_3.xlsx edit:
# - All metric rows (B:Q, rows 2-26) get overwritten with a new, identical set of
#   metric values (the "best mse_vt" model run).
# - Column A (model name) gets reordered according to the new ranking, except
#   row 13 which keeps its original "model_4_3_23" label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order of model names for rows 2..26
$newOrder = @(
    "model_4_3_0",
    "model_4_3_22",
    "model_4_3_21",
    "model_4_3_20",
    "model_4_3_19",
    "model_4_3_18",
    "model_4_3_17",
    "model_4_3_16",
    "model_4_3_15",
    "model_4_3_14",
    "model_4_3_13",
    "model_4_3_23",
    "model_4_3_12",
    "model_4_3_10",
    "model_4_3_9",
    "model_4_3_8",
    "model_4_3_7",
    "model_4_3_6",
    "model_4_3_5",
    "model_4_3_4",
    "model_4_3_3",
    "model_4_3_2",
    "model_4_3_1",
    "model_4_3_11",
    "model_4_3_24"
)

# New metric values (columns B..Q), identical for every row
$newValues = @(
    0.9994912738785552,
    0.6806725964963172,
    0.9993476784002224,
    0.9962375791921385,
    0.9986413900465918,
    0.0003020018015138525,
    0.1895665409453636,
    0.0007003396042700758,
    0.001187171653545526,
    0.000943755628907801,
    0.007272400149046574,
    0.01737819902964207,
    1.000329984511207,
    0.01811802494016521,
    138.2101551506591,
    212.5615804676193
)

for ($i = 0; $i -lt $newOrder.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newOrder[$i]

    for ($j = 0; $j -lt $newValues.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $newValues[$j]
    }
}
